# Update countries & provincias Spain
# This script applies the data refresh captured in the commit:
#  - Update the "last updated" timestamp string
#  - Refresh a handful of per-country statistics (Estados Unidos totals,
#    Monaco, and a couple of small territories)
#  - Reorder two pairs of small territories (Santa Lucia/Belice and
#    Montserrat/Groenlandia), which also carries their own stats with them

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 20:35"

# 2) Update Estados Unidos totals (row 4)
$ws.Cells.Item(4, 2).Value = 1579137   # Casos totales
$ws.Cells.Item(4, 3).Value = 8554      # Nuevos casos
$ws.Cells.Item(4, 5).Value = 1120309   # Recuperados
$ws.Cells.Item(4, 7).Value = 612       # Muertes hoy
$ws.Cells.Item(4, 8).Value = 94145     # Muertes

# 3) Update Monaco totals (row 170)
$ws.Cells.Item(170, 4).Value = 90      # Casos activos
$ws.Cells.Item(170, 5).Value = 3       # Recuperados

# 4) Swap Santa Lucia (row 196) and Belice (row 197), keeping each
#    country's own stats attached to its name
$ws.Cells.Item(196, 1).Value = "Belice"
$ws.Cells.Item(196, 4).Value = 16
$ws.Cells.Item(196, 8).Value = 2

$ws.Cells.Item(197, 1).Value = "Santa Lucia"
$ws.Cells.Item(197, 4).Value = 18
$ws.Cells.Item(197, 8).Value = 0

# 5) Swap Montserrat (row 209) and Groenlandia (row 210), keeping each
#    country's own stats attached to its name
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 0

$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1
